$wb = $excel.ActiveWorkbook

# ===========================================================================
# Two new files were handed off:
#   e0a98168-8dc6-480c-96e4-feb5bafc2ad5.md
#   ede29265-f53d-414c-9b14-52b38b8de1ae.md
# This adds a row for each to the "Overview" sheet and to each locale sheet
# (zh-cn, de-de), expanding every table from 5 to 7 data rows.
# ===========================================================================

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows("5").Copy()
$wsOverview.Rows("6").Insert()
$wsOverview.Rows("5").Copy()
$wsOverview.Rows("7").Insert()

$wsOverview.Range("A6").Value = "e0a98168-8dc6-480c-96e4-feb5bafc2ad5.md"
$wsOverview.Range("B6").Value = "e2e\e0a98168-8dc6-480c-96e4-feb5bafc2ad5.md"
$wsOverview.Range("G6").Value = "2016-09-04 20:47:20"

$wsOverview.Range("A7").Value = "ede29265-f53d-414c-9b14-52b38b8de1ae.md"
$wsOverview.Range("B7").Value = "e2e\ede29265-f53d-414c-9b14-52b38b8de1ae.md"
$wsOverview.Range("G7").Value = "2016-09-04 20:47:20"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e0a98168/e2e/e0a98168-8dc6-480c-96e4-feb5bafc2ad5.md", "", "", "e2e\e0a98168-8dc6-480c-96e4-feb5bafc2ad5.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ede29265/e2e/ede29265-f53d-414c-9b14-52b38b8de1ae.md", "", "", "e2e\ede29265-f53d-414c-9b14-52b38b8de1ae.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G7"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Rows("5").Copy()
$wsZhCn.Rows("6").Insert()
$wsZhCn.Rows("5").Copy()
$wsZhCn.Rows("7").Insert()

$wsZhCn.Range("A6").Value = "e0a98168-8dc6-480c-96e4-feb5bafc2ad5.md"
$wsZhCn.Range("G6").Value = "e0a98168-8dc6-480c-96e4-feb5bafc2ad5.1cd7158fa032a221835d69836ba5a68086bf0141.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-09-04 20:47:16"

$wsZhCn.Range("A7").Value = "ede29265-f53d-414c-9b14-52b38b8de1ae.md"
$wsZhCn.Range("G7").Value = "ede29265-f53d-414c-9b14-52b38b8de1ae.a9e0385c7487949634b4ae78cfedd6f108c9b6c9.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-09-04 20:47:16"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e0a98168/e2e/e0a98168-8dc6-480c-96e4-feb5bafc2ad5.md", "", "", "e0a98168-8dc6-480c-96e4-feb5bafc2ad5.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ede29265/e2e/ede29265-f53d-414c-9b14-52b38b8de1ae.md", "", "", "ede29265-f53d-414c-9b14-52b38b8de1ae.md")

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P7"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Rows("5").Copy()
$wsDeDe.Rows("6").Insert()
$wsDeDe.Rows("5").Copy()
$wsDeDe.Rows("7").Insert()

$wsDeDe.Range("A6").Value = "e0a98168-8dc6-480c-96e4-feb5bafc2ad5.md"
$wsDeDe.Range("G6").Value = "e0a98168-8dc6-480c-96e4-feb5bafc2ad5.1cd7158fa032a221835d69836ba5a68086bf0141.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-09-04 20:47:20"

$wsDeDe.Range("A7").Value = "ede29265-f53d-414c-9b14-52b38b8de1ae.md"
$wsDeDe.Range("G7").Value = "ede29265-f53d-414c-9b14-52b38b8de1ae.a9e0385c7487949634b4ae78cfedd6f108c9b6c9.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-09-04 20:47:20"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e0a98168/e2e/e0a98168-8dc6-480c-96e4-feb5bafc2ad5.md", "", "", "e0a98168-8dc6-480c-96e4-feb5bafc2ad5.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ede29265/e2e/ede29265-f53d-414c-9b14-52b38b8de1ae.md", "", "", "ede29265-f53d-414c-9b14-52b38b8de1ae.md")

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P7"))

Write-Output "Handoff report rows added to Overview, zh-cn, de-de"
